$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was added to the top of the data table (row 13),
# pushing all the existing data rows down by one (old row 13 -> 14, ...,
# old row 73 -> 74).
$ws.Rows.Item(13).Insert()

$ws.Range("A13").Value = 7
$ws.Range("B13").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C13").Value = "Ñuble"
$ws.Range("D13").Value = 44623
$ws.Range("E13").Value = 16
$ws.Range("F13").Value = 100112031
$ws.Range("G13").Value = "Poroto verde"
$ws.Range("H13").Value = "Sin especificar"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 60
$ws.Range("K13").Value = 28000
$ws.Range("L13").Value = 29000
$ws.Range("M13").Value = 28500
$ws.Range("N13").Value = "$/saco 25 kilos"
$ws.Range("O13").Value = "Región del Maule"
$ws.Range("P13").Value = 1140
$ws.Range("Q13").Value = 25
$ws.Range("R13").Value = "Hortaliza"
